# Natmi following Dr Hou advice
# Expand the Agrn-Atp1a3 LR-pair sheet from "only matching sending/target
# cluster" rows to the full 3x3 cross join of clusters (ECs, FAPs, sCs) for
# sending cluster (A) x target cluster (D), recomputing all derived columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Agrn", "Atp1a3", "ECs",  3, 1, 15.59657533333333, 46.789726,   0.4757744772251148,  0.475774477225115,   2, 0.6666666666666666, 11.51598433333333,   34.547953,  0.9063460549271956,  0.9063460549271956,  179.6099171923198, 1616.489254730878, 0.4312163204680317,  0.4312163204680318),
    @("ECs",  "Agrn", "Atp1a3", "FAPs", 3, 1, 15.59657533333333, 46.789726,   0.4757744772251148,  0.475774477225115,   3, 1,                  1.181969666666667,   3.545909,   0.09302492200567822, 0.09302492200567822, 18.43467894788155, 165.912110530934,  0.04425888363615864, 0.04425888363615865),
    @("ECs",  "Agrn", "Atp1a3", "sCs",  3, 1, 15.59657533333333, 46.789726,   0.4757744772251148,  0.475774477225115,   2, 0.6666666666666666, 0.007992333333333334, 0.023977, 0.0006290230671261296, 0.0006290230671261296, 0.1246530289224445, 1.121877260302, 0.0002992731209244726, 0.0002992731209244727),
    @("FAPs", "Agrn", "Atp1a3", "ECs",  3, 1, 6.399531333333333, 19.198594,   0.1952180917624358,  0.1952180917624358,  2, 0.6666666666666666, 11.51598433333333,   34.547953,  0.9063460549271956,  0.9063460549271956,  73.69690257534245, 663.272123178082,  0.1769351473192989,  0.1769351473192989),
    @("FAPs", "Agrn", "Atp1a3", "FAPs", 3, 1, 6.399531333333333, 19.198594,   0.1952180917624358,  0.1952180917624358,  3, 1,                  1.181969666666667,   3.545909,   0.09302492200567822, 0.09302492200567822, 7.564051916882888, 68.076467251946,   0.01816014776029792, 0.01816014776029792),
    @("FAPs", "Agrn", "Atp1a3", "sCs",  3, 1, 6.399531333333333, 19.198594,   0.1952180917624358,  0.1952180917624358,  2, 0.6666666666666666, 0.007992333333333334, 0.023977, 0.0006290230671261296, 0.0006290230671261296, 0.05114718759311111, 0.4603246883380001, 0.0001227966828389176, 0.0001227966828389176),
    @("sCs",  "Agrn", "Atp1a3", "ECs",  3, 1, 10.78533933333333, 32.356018,   0.3290074310124493,  0.3290074310124493,  2, 0.6666666666666666, 11.51598433333333,   34.547953,  0.9063460549271956,  0.9063460549271956,  124.2037987923504, 1117.834189131154, 0.2981945871398649,  0.2981945871398649),
    @("sCs",  "Agrn", "Atp1a3", "FAPs", 3, 1, 10.78533933333333, 32.356018,   0.3290074310124493,  0.3290074310124493,  3, 1,                  1.181969666666667,   3.545909,   0.09302492200567822, 0.09302492200567822, 12.74794393670689, 114.731495430362,  0.03060589060922165, 0.03060589060922166),
    @("sCs",  "Agrn", "Atp1a3", "sCs",  3, 1, 10.78533933333333, 32.356018,   0.3290074310124493,  0.3290074310124493,  2, 0.6666666666666666, 0.007992333333333334, 0.023977, 0.0006290230671261296, 0.0006290230671261296, 0.08620002706511111, 0.7758002435860001, 0.0002069532633627393, 0.0002069532633627394)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
